# Auto-generated Excel COM-interop script
# Applies the numeric corrections described by the commit diff to the
# "Zalera_Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 14297296
$ws.Range("J76").Value = 11894
$ws.Range("L76").Value = 11894
$ws.Range("N76").Value = -12524

$ws.Range("H79").Value = 14297296
$ws.Range("J79").Value = 11894
$ws.Range("L79").Value = 11894
$ws.Range("N79").Value = -14078

$ws.Range("H98").Value = 7036.0415
$ws.Range("I98").Value = 8376.632
$ws.Range("K98").Value = 8376.632
$ws.Range("M98").Value = -6878.632

$ws.Range("H101").Value = 3088.2222
$ws.Range("I101").Value = 3257
$ws.Range("J101").Value = 2497.5
$ws.Range("K101").Value = 9771
$ws.Range("L101").Value = 7492.5
$ws.Range("M101").Value = -8149
$ws.Range("N101").Value = -10736.5

$ws.Range("H107").Value = 25001322
$ws.Range("I107").Value = 25001322
$ws.Range("K107").Value = 25001322
$ws.Range("M107").Value = -24999402

$ws.Range("H122").Value = 7036.0415
$ws.Range("I122").Value = 8376.632
$ws.Range("K122").Value = 25129.896
$ws.Range("M122").Value = -22679.896

$ws.Range("H129").Value = 1843
$ws.Range("I129").Value = 1076
$ws.Range("J129").Value = 2695.2222
$ws.Range("K129").Value = 3228
$ws.Range("L129").Value = 8085.6666
$ws.Range("M129").Value = 1772
$ws.Range("N129").Value = -18085.6666

$ws.Range("H132").Value = 3895.8333
$ws.Range("I132").Value = 2559
$ws.Range("J132").Value = 5232.6665
$ws.Range("K132").Value = 7677
$ws.Range("L132").Value = 15697.9995
$ws.Range("M132").Value = -5147
$ws.Range("N132").Value = -20757.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 42949.332
$ws.Range("I32").Value = 61199.57
$ws.Range("K32").Value = 61199.57
$ws.Range("M32").Value = -60912.57

$ws.Range("H61").Value = 12350674
$ws.Range("I61").Value = 15154869
$ws.Range("K61").Value = 15154869
$ws.Range("M61").Value = -15154657

$ws.Range("H74").Value = 388748.03
$ws.Range("I74").Value = 528067.7
$ws.Range("K74").Value = 528067.7
$ws.Range("M74").Value = -527193.7

$ws.Range("H77").Value = 388748.03
$ws.Range("I77").Value = 528067.7
$ws.Range("K77").Value = 2640338.5
$ws.Range("M77").Value = -2635970.5

$ws.Range("H132").Value = 4005410.5
$ws.Range("I132").Value = 5409097.5
$ws.Range("K132").Value = 16227292.5
$ws.Range("M132").Value = -16224762.5

$ws.Range("H136").Value = 12350674
$ws.Range("I136").Value = 15154869
$ws.Range("K136").Value = 45464607
$ws.Range("M136").Value = -45462057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 26666.666
$ws.Range("I75").Value = 26666.666
$ws.Range("K75").Value = 26666.666
$ws.Range("M75").Value = -25730.666

$ws.Range("H76").Value = 21314
$ws.Range("J76").Value = 21314
$ws.Range("L76").Value = 21314
$ws.Range("N76").Value = -21944

$ws.Range("H78").Value = 26666.666
$ws.Range("I78").Value = 26666.666
$ws.Range("K78").Value = 79999.99800000001
$ws.Range("M78").Value = -75319.99800000001

$ws.Range("H79").Value = 21314
$ws.Range("J79").Value = 21314
$ws.Range("L79").Value = 21314
$ws.Range("N79").Value = -23498

$ws.Range("H88").Value = 33495.234
$ws.Range("J88").Value = 33495.234
$ws.Range("L88").Value = 33495.234
$ws.Range("N88").Value = -34307.234

$ws.Range("H91").Value = 33495.234
$ws.Range("J91").Value = 33495.234
$ws.Range("L91").Value = 33495.234
$ws.Range("N91").Value = -36303.234

$ws.Range("H134").Value = 8285.883
$ws.Range("I134").Value = 7190.6665
$ws.Range("K134").Value = 21571.9995
$ws.Range("M134").Value = -19036.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 56380.332
$ws.Range("I12").Value = 62802.875
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 62802.875
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -62632.875
$ws.Range("N12").Value = -5340

$ws.Range("H31").Value = 5541.607
$ws.Range("I31").Value = 3399.111
$ws.Range("J31").Value = 6556.4736
$ws.Range("K31").Value = 3399.111
$ws.Range("L31").Value = 6556.4736
$ws.Range("M31").Value = -3104.111
$ws.Range("N31").Value = -7146.4736

$ws.Range("H34").Value = 5541.607
$ws.Range("I34").Value = 3399.111
$ws.Range("J34").Value = 6556.4736
$ws.Range("K34").Value = 3399.111
$ws.Range("L34").Value = 6556.4736
$ws.Range("M34").Value = -3197.111
$ws.Range("N34").Value = -6960.4736

$ws.Range("H58").Value = 3865
$ws.Range("J58").Value = 5787.4546
$ws.Range("L58").Value = 5787.4546
$ws.Range("N58").Value = -6193.4546

$ws.Range("H122").Value = 2373.5
$ws.Range("J122").Value = 4284.5713
$ws.Range("L122").Value = 12853.7139
$ws.Range("N122").Value = -17753.7139

$ws.Range("H132").Value = 28848.875
$ws.Range("I132").Value = 5467.0586
$ws.Range("J132").Value = 85633.28999999999
$ws.Range("K132").Value = 16401.1758
$ws.Range("L132").Value = 256899.87
$ws.Range("M132").Value = -13871.1758
$ws.Range("N132").Value = -261959.87

$ws.Range("H134").Value = 6833.2915
$ws.Range("I134").Value = 4742.3076
$ws.Range("K134").Value = 14226.9228
$ws.Range("M134").Value = -11691.9228

$ws.Range("H136").Value = 3865
$ws.Range("J136").Value = 5787.4546
$ws.Range("L136").Value = 17362.3638
$ws.Range("N136").Value = -22462.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 409.23077
$ws.Range("I23").Value = 91.125
$ws.Range("J23").Value = 918.2
$ws.Range("K23").Value = 273.375
$ws.Range("L23").Value = 2754.6
$ws.Range("M23").Value = -38.375
$ws.Range("N23").Value = -3224.6

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H86").Value = 2338.25
$ws.Range("J86").Value = 4350
$ws.Range("L86").Value = 13050
$ws.Range("N86").Value = -15422

$ws.Range("H89").Value = 2338.25
$ws.Range("J89").Value = 4350
$ws.Range("L89").Value = 39150
$ws.Range("N89").Value = -51006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 6809
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 6809
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 6809
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -7035

$ws.Range("H16").Value = 6809
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 6809
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 6809
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -7309

$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16622

$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -53112

$ws.Range("H132").Value = 7402.9414
$ws.Range("I132").Value = 4526.923
$ws.Range("J132").Value = 16750
$ws.Range("K132").Value = 13580.769
$ws.Range("L132").Value = 50250
$ws.Range("M132").Value = -11050.769
$ws.Range("N132").Value = -55310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4223.75
$ws.Range("I10").Value = 998.3333
$ws.Range("K10").Value = 998.3333
$ws.Range("M10").Value = -858.3333

$ws.Range("H122").Value = 4034.7144
$ws.Range("I122").Value = 4049.2
$ws.Range("J122").Value = 3998.5
$ws.Range("K122").Value = 12147.6
$ws.Range("L122").Value = 11995.5
$ws.Range("M122").Value = -9697.599999999999
$ws.Range("N122").Value = -16895.5

$ws.Range("H136").Value = 3180071.5
$ws.Range("I136").Value = 5132550
$ws.Range("K136").Value = 15397650
$ws.Range("M136").Value = -15395100

$ws.Range("H140").Value = 99998.664
$ws.Range("I140").Value = 90000
$ws.Range("J140").Value = 101998.4
$ws.Range("K140").Value = 90000
$ws.Range("L140").Value = 101998.4
$ws.Range("M140").Value = -84820
$ws.Range("N140").Value = -112358.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 72000
$ws.Range("J110").Value = 72000
$ws.Range("L110").Value = 72000
$ws.Range("N110").Value = -80180

$ws.Range("H136").Value = 3664963.2
$ws.Range("I136").Value = 4330329.5
$ws.Range("K136").Value = 12990988.5
$ws.Range("M136").Value = -12988438.5

